$p = $ppt.ActivePresentation

# Remove the last 9 slides (slides 14 through 22) from the presentation.
for ($i = $p.Slides.Count; $i -ge 14; $i--) {
    $p.Slides.Item($i).Delete()
}
